$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.215.06"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "4.004.42"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'596.90"
$ws.Range("E5").Value = "  +10.46%  "
$ws.Range("D6").Value = "'160.42"
$ws.Range("E6").Value = "  +7.23%  "
$ws.Range("D7").Value = "'0.685"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.752"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").Value = "'54.10"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "'11.01"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "4.648.00"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "3.994.81"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "'1.28"
$ws.Range("E16").Value = "  +9.09%  "
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "'20.40"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "72.952.77"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "'435.33"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "  +12.06%  "
$ws.Range("D23").Value = "'96.39"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'4.41"
$ws.Range("E26").Value = "  +16.67%  "
$ws.Range("D27").Value = "'11.24"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'10.48"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").Value = "'36.50"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "'13.84"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "'0.130"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "'48.39"
$ws.Range("E34").Value = "  -4.93%  "
$ws.Range("D35").Value = "'671.15"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").Value = "'70.75"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("D37").Value = "0.0₃0911"
$ws.Range("E37").Value = "  +11.05%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "'3.35"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0494"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'10.68"
$ws.Range("E45").Value = "  +9.82%  "
$ws.Range("D46").Value = "'0.150"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "'2.61"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").Value = "'3.39"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "2.874.76"
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "'3.40"
$ws.Range("E51").Value = "  +4.48%  "
